# Weekly data refresh: re-pull of "Fruta / hortaliza" series shuffles the
# existing 14 data rows (each keeps its original record, but rows are
# re-ordered/re-dated as the upstream weekly snapshot rotates) and appends a
# new observation as row 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column map: D=Fecha, J=Volumen, K=Precio minimo, L=Precio maximo,
#             M=Precio promedio ponderado, P=Precio $/Kg

# row -> @(Fecha, Volumen, PrecioMin, PrecioProm, PrecioKg)
# (L / Precio maximo is unchanged for every existing row, so it is left alone)
$updates = @{
    2  = @(44350, 25,  6000, 6000, 375)
    3  = @(44306, 50,  6000, 6000, 375)
    4  = @(44403, 43,  6000, 6000, 375)
    6  = @(44363, 160, 5500, 5750, 359)
    7  = @(44355, 25,  6000, 6000, 375)
    8  = @(44407, 45,  5500, 5744, 359)
    9  = @(44330, 120, 6000, 6000, 375)
    10 = @(44328, 160, 6000, 6000, 375)
    11 = @(44371, 34,  5500, 5750, 359)
    12 = @(44358, 52,  6000, 6000, 375)
    13 = @(44341, 51,  5500, 5755, 360)
    14 = @(44313, 34,  6000, 6000, 375)
    15 = @(44438, 34,  5000, 5500, 344)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D: Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]   # K: Precio minimo
    $ws.Cells.Item($row, 13).Value = $vals[3]   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals[4]   # P: Precio $/Kg
}

# New weekly observation appended as row 16
$ws.Cells.Item(16, 1).Value  = 9
$ws.Cells.Item(16, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(16, 3).Value  = "Metropolitana"
$ws.Cells.Item(16, 4).Value  = 44442
$ws.Cells.Item(16, 5).Value  = 13
$ws.Cells.Item(16, 6).Value  = 100112010
$ws.Cells.Item(16, 7).Value  = "Achicoria"
$ws.Cells.Item(16, 8).Value  = "Sin especificar"
$ws.Cells.Item(16, 9).Value  = "Primera"
$ws.Cells.Item(16, 10).Value = 25
$ws.Cells.Item(16, 11).Value = 6000
$ws.Cells.Item(16, 12).Value = 7000
$ws.Cells.Item(16, 13).Value = 6480
$ws.Cells.Item(16, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(16, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(16, 16).Value = 405
$ws.Cells.Item(16, 17).Value = 16
$ws.Cells.Item(16, 18).Value = "Hortaliza"

# Match the date-column number format used by the rest of the Fecha column
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(15, 4).NumberFormat
